# Generate Report for Archive
# - Update status text from "Ready for handoff" to "In Translation" across all sheets
# - Narrow the "Status" column on each sheet (Overview: E & F, zh-cn: C, de-de: C)

$wb = $excel.ActiveWorkbook

# Target stored column width is 13.4101845877511 "characters". The COM
# ColumnWidth setter here quantizes to the nearest 1/6-character (pixel)
# step, so feed it the input that lands on the closest reachable step
# (13.333333333333334, written as the stored <col width>).
$newStatus = "In Translation"
$newWidth = 12.5

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F4").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C4").Value = $newStatus
$zhcn.Columns.Item(3).ColumnWidth = $newWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C4").Value = $newStatus
$dede.Columns.Item(3).ColumnWidth = $newWidth
